{"js": "// UC1 fix: the use-case step text had been split across three runs with a\n// stray \"(\" run in the middle, e.g.:\n//   \"...Acquista biglietto sessione\"  +  \"(\"  +  \"\".\n// The intended text has no stray parenthesis:\n//   \"...Acquista biglietto sessione\".\n//\n// Find the exact broken fragment and replace it with the corrected text.\nconst body = context.document.body;\n\nconst searchResults = body.search(\"sessione(\u201d.\", { matchCase: true, matchWildcards: false });\nsearchResults.load(\"items,text\");\nawait context.sync();\n\nif (searchResults.items.length !== 1) {\n  throw new Error(\n    `Expected exactly one match for the broken \"sessione(\u201d.\" fragment, found ${searchResults.items.length}.`\n  );\n}\n\nsearchResults.items[0].insertText(\"sessione\u201d.\", \"Replace\");\nawait context.sync();\n", "ps1": "# UC1 fix: the use-case step text had been split across three runs with a\n# stray \"(\" run in the middle, e.g.:\n#   \"...Acquista biglietto sessione\"  +  \"(\"  +  \"\".\n# The intended text has no stray parenthesis:\n#   \"...Acquista biglietto sessione\".\n#\n# Use Find/Replace on the document content to correct the broken fragment.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"sessione(\u201d.\"\n$find.Replacement.Text = \"sessione\u201d.\"\n$find.Forward = $true\n$find.Wrap = 0  # wdFindStop - there should be exactly one occurrence\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 0, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw \"Could not find the broken 'sessione(\u201d.' fragment to fix.\"\n}\n"}
